$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B-column renumbering (two node-ID values were bumped) ---
$ws.Range("B17").Value = 4471
$ws.Range("B117").Value = 5471

# --- Newly populated rows: Old ID (A) / New Label (C) for previously-blank rows ---
$rows = @(
    @{ Row = 32; A = 7;   C = "Henri Morriseau" },
    @{ Row = 33; A = 278; C = "Lafayette St. Pierre" },
    @{ Row = 34; A = 349; C = "Doctor's Daugther" },
    @{ Row = 35; A = 350; C = "Local Doctor" },
    @{ Row = 42; A = 51;  C = "Stan Anderson" },
    @{ Row = 43; A = 80;  C = "Stan's Father" },
    @{ Row = 44; A = 81;  C = "Stan's Mother" },
    @{ Row = 45; A = 5;   C = "Narrator" },
    @{ Row = 46; A = 53;  C = "Narrator's Father" },
    @{ Row = 47; A = 55;  C = "Narrator's Grandfather" },
    @{ Row = 48; A = 56;  C = "Narrator's Grandmother" },
    @{ Row = 49; A = 52;  C = "Narrator's Mother" },
    @{ Row = 50; A = 28;  C = "Stan's Brother" },
    @{ Row = 51; C = "Stan's Sister" },
    @{ Row = 62; A = 19;  C = "Shamengwa's Daughter" },
    @{ Row = 63; A = 7;   C = "Baptiste Parentheau" },
    @{ Row = 64; A = 330; C = "Baptiste Parentheau's Father" },
    @{ Row = 65; A = 333; C = "Baptiste Parentheau's Uncles" },
    @{ Row = 66; A = 331; C = "Baptiste Parentheau's Wife" },
    @{ Row = 67; A = 278; C = "Edwin Parentheau" },
    @{ Row = 68; A = 332; C = "Edwin Parentheau's Wife" },
    @{ Row = 69; A = 4;   C = "Judge (Narrator)" },
    @{ Row = 70; A = 96;  C = "Shamengwa's Father" },
    @{ Row = 71; A = 95;  C = "Shamengwa's Mother" },
    @{ Row = 72; A = 101; C = "Baby who was lost" },
    @{ Row = 73; A = 15;  C = "Shamengwa's Sister" },
    @{ Row = 74; C = "Corwin's roommates" },
    @{ Row = 75; C = "Priest" },
    @{ Row = 82; A = 334; C = "Sprecht " },
    @{ Row = 83; A = 4;   C = "Narrator" },
    @{ Row = 84; A = 250; C = "C.'s Patients" },
    @{ Row = 85; A = 338; C = "C.'s Receptionist" },
    @{ Row = 86; A = 6;   C = "C." },
    @{ Row = 87; A = 281; C = "Narrator's Mother" },
    @{ Row = 88; A = 225; C = "Narrator's Father" }
)

foreach ($r in $rows) {
    if ($r.ContainsKey("A")) {
        $ws.Cells.Item($r.Row, 1).Value = $r.A
    }
    $ws.Cells.Item($r.Row, 3).Value = $r.C
}

# Row 82's label ("Sprecht ") is set in a different font (Garamond) than the rest of the sheet.
$ws.Range("C82").Font.Name = "Garamond"

# --- View state: restore the scroll position / selection the author left the sheet in ---
$ws.Range("A62").Select()
$ws.Range("F83").Select()
